$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell address -> new literal text value.
# We route every write through Formula (string literal) + Copy + PasteSpecial(Values)
# so the result lands as a plain text value (matching the original inlineStr cells)
# instead of letting Excels smart entry reinterpret numeric-looking text as a number
# (which would both lose exact formatting like trailing zeros and silently add a
# NumberFormat style that is not present in the source workbook).
$updates = @(
    @{Cell = 'D2'; Value = '37.371.72'}
    @{Cell = 'E2'; Value = '  -0.08%  '}
    @{Cell = 'D3'; Value = '2.067.68'}
    @{Cell = 'E3'; Value = '  +0.03%  '}
    @{Cell = 'D5'; Value = '234.88'}
    @{Cell = 'E5'; Value = '  -0.33%  '}
    @{Cell = 'D6'; Value = '0.625'}
    @{Cell = 'E6'; Value = '  +1.90%  '}
    @{Cell = 'E7'; Value = '  +0.04%  '}
    @{Cell = 'D8'; Value = '57.22'}
    @{Cell = 'E8'; Value = '  -1.97%  '}
    @{Cell = 'E9'; Value = '  +2.81%  '}
    @{Cell = 'D10'; Value = '0.0772'}
    @{Cell = 'E10'; Value = '  +1.20%  '}
    @{Cell = 'D11'; Value = '0.103'}
    @{Cell = 'E11'; Value = '  +0.88%  '}
    @{Cell = 'D12'; Value = '2.370.83'}
    @{Cell = 'E12'; Value = '  +0.02%  '}
    @{Cell = 'D13'; Value = '14.39'}
    @{Cell = 'E13'; Value = '  -1.62%  '}
    @{Cell = 'D14'; Value = '20.61'}
    @{Cell = 'E14'; Value = '  -1.67%  '}
    @{Cell = 'E15'; Value = '  -0.65%  '}
    @{Cell = 'D16'; Value = '5.19'}
    @{Cell = 'E16'; Value = '  -0.45%  '}
    @{Cell = 'D17'; Value = '2.067.97'}
    @{Cell = 'E17'; Value = '  -0.03%  '}
    @{Cell = 'D18'; Value = '37.316.21'}
    @{Cell = 'E18'; Value = '  -0.78%  '}
    @{Cell = 'E19'; Value = '  -0.33%  '}
    @{Cell = 'D20'; Value = '69.55'}
    @{Cell = 'E20'; Value = '  +0.73%  '}
    @{Cell = 'E21'; Value = '  +0.08%  '}
    @{Cell = 'D22'; Value = '226.73'}
    @{Cell = 'E22'; Value = '  +0.02%  '}
    @{Cell = 'E23'; Value = '  +0.07%  '}
    @{Cell = 'E24'; Value = '  +1.73%  '}
    @{Cell = 'D25'; Value = '2.39'}
    @{Cell = 'E25'; Value = '  -2.67%  '}
    @{Cell = 'D26'; Value = '167.00'}
    @{Cell = 'E26'; Value = '  +1.45%  '}
    @{Cell = 'D27'; Value = '8.92'}
    @{Cell = 'E27'; Value = '  +0.42%  '}
    @{Cell = 'D28'; Value = '1.40'}
    @{Cell = 'E28'; Value = '  -6.63%  '}
    @{Cell = 'D29'; Value = '0.129'}
    @{Cell = 'E29'; Value = '  +0.87%  '}
    @{Cell = 'D30'; Value = '19.11'}
    @{Cell = 'E30'; Value = '  -0.70%  '}
    @{Cell = 'D31'; Value = '0.118'}
    @{Cell = 'E31'; Value = '  -1.06%  '}
    @{Cell = 'D32'; Value = '4.51'}
    @{Cell = 'E32'; Value = '  +0.30%  '}
    @{Cell = 'E33'; Value = '  -1.21%  '}
    @{Cell = 'D34'; Value = '4.52'}
    @{Cell = 'E34'; Value = '  +0.24%  '}
    @{Cell = 'E35'; Value = '  -3.34%  '}
    @{Cell = 'E36'; Value = '  +0.43%  '}
    @{Cell = 'E37'; Value = '  -2.99%  '}
    @{Cell = 'E38'; Value = '  -0.04%  '}
    @{Cell = 'D39'; Value = '5.60'}
    @{Cell = 'E39'; Value = '  -5.15%  '}
    @{Cell = 'D40'; Value = '2.95'}
    @{Cell = 'E40'; Value = '  -0.50%  '}
    @{Cell = 'D41'; Value = '0.0956'}
    @{Cell = 'E41'; Value = '  -2.84%  '}
    @{Cell = 'B42'; Value = 'Maker'}
    @{Cell = 'C42'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'}
    @{Cell = 'D42'; Value = '1.486.48'}
    @{Cell = 'E42'; Value = '  +0.49%  '}
    @{Cell = 'B43'; Value = 'Aave'}
    @{Cell = 'C43'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'}
    @{Cell = 'D43'; Value = '97.50'}
    @{Cell = 'E43'; Value = '  +0.72%  '}
    @{Cell = 'E44'; Value = '  +0.87%  '}
    @{Cell = 'E45'; Value = '  -0.86%  '}
    @{Cell = 'D46'; Value = '4.14'}
    @{Cell = 'E46'; Value = '  -8.03%  '}
    @{Cell = 'E47'; Value = '  -0.67%  '}
    @{Cell = 'D48'; Value = '7.19'}
    @{Cell = 'E48'; Value = '  -1.53%  '}
    @{Cell = 'D49'; Value = '15.04'}
    @{Cell = 'E49'; Value = '  -5.44%  '}
    @{Cell = 'E50'; Value = '  +0.82%  '}
    @{Cell = 'D51'; Value = '47.49'}
    @{Cell = 'E51'; Value = '  +5.88%  '}
)

foreach ($u in $updates) {
    $target = $ws.Range($u.Cell)
    $escaped = $u.Value.Replace('"', '""')
    $target.Formula = '="' + $escaped + '"'
    $target.Copy()
    $target.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false

